# Chronologie2023.xlsx — "Add files via upload"
#
# Two new frog/toad survey observations (217 & 227 route, Saint-Félicien,
# Saguenay - Lac-Saint-Jean, zone D, Cote 3, contact Alexandra Coutu, date
# 2023-04-28) are added to the bottom of the data table, then the whole
# table is re-sorted by Date then Route associée (as it apparently was
# before), which re-shuffles the previously-last rows further down.
# Finally the view selection is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the two new rows at the bottom of the table (rows 39-40) ---

# Borrow the formatting of two existing rows that already carry the PSCR /
# RASY species colouring used elsewhere in the sheet, then fix up the
# "Zone climatique" (column F) formatting to the "D" zone colour.
$ws.Range("A17:I17").Copy()
$ws.Range("A39:I39").PasteSpecial(-4122)
$ws.Range("A18:I18").Copy()
$ws.Range("A40:I40").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("F40").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 39: PSCR observation
$ws.Range("A39").Value2 = 45044
$ws.Range("B39").Value = "PSCR"
$ws.Range("C39").Value = "217 & 227"
$ws.Range("D39").Value = "Saint-Félicien"
$ws.Range("E39").Value = "Saguenay - Lac-Saint-Jean"
$ws.Range("F39").Value = "D"
$ws.Range("G39").Value = "Cote 3"
$ws.Range("H39").Value = ""
$ws.Range("I39").Value = "Alexandra Coutu"

# Row 40: RASY observation (same route/city/date/contact)
$ws.Range("A40").Value2 = 45044
$ws.Range("B40").Value = "RASY"
$ws.Range("C40").Value = "217 & 227"
$ws.Range("D40").Value = "Saint-Félicien"
$ws.Range("E40").Value = "Saguenay - Lac-Saint-Jean"
$ws.Range("F40").Value = "D"
$ws.Range("G40").Value = "Cote 3"
$ws.Range("H40").Value = ""
$ws.Range("I40").Value = "Alexandra Coutu"

# --- 2. Re-sort the whole data range by Date (A) then Route associée (F) ---

$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A4:A40"))
$sort.SortFields.Add($ws.Range("F4:F40"))
$sort.SetRange($ws.Range("A4:I40"))
$sort.Header = -4142
$sort.Apply()

# --- 3. Update the active selection on the sheet ---

$null = $ws.Range("N26").Select()
